$d = $word.ActiveDocument

# ---- Clear all existing content down to a single trailing paragraph ----
while ($d.Paragraphs.Count -gt 1) {
    $p = $d.Paragraphs(1)
    $full = $d.Range($p.Range.Start, $p.Range.End)
    $full.Delete()
}

# The single remaining paragraph still carries the original paragraph's
# rsid/paraId attributes. Insert a brand-new paragraph after it and drop
# the old one so every paragraph we build from here on is attribute-free,
# matching freshly-authored content.
$d.Paragraphs(1).Range.InsertParagraphAfter()
$oldFirst = $d.Paragraphs(1)
$d.Range($oldFirst.Range.Start, $oldFirst.Range.End).Delete()

# ---- Paragraph text content (in order) ----
$texts = @(
    "Ingredients\n",
    "`tGarlic\n",
    "\n",
    "Instructions\n",
    "Drying garlic is an easy and safe way to keep garlic indefinitely. Choose fresh, firm, flavourful cloves with no bruises and follow these simple directions.\n",
    "\n",
    "Peel the garlic cloves.\n",
    "\n",
    "Slice the garlic thinly. A food processor works well.\n",
    "\n",
    "Dry the garlic until crisp. In a dehydrator you have a choice of faster, hotter drying or slower, cooler drying. We dry our garlic for two days at about 45°C (115°F).\n",
    "\n",
    "Dried garlic may be stored at room temperature in an airtight container.\n",
    "\n",
    "If you store the dried garlic in the freezer in the form of flakes, and then grind it close to the time when you will be using the garlic it will keep its amazing freshness for more than a year.\n",
    "\n "
)

# Paragraphs 1-4 (Ingredients / Garlic / blank / Instructions) have no left indent.
# Paragraphs 5-16 (the instruction lines) are indented 720 twips (36pt) to the left.
$indented = @($false, $false, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true)

for ($i = 0; $i -lt $texts.Count; $i++) {
    $idx = $i + 1
    if ($idx -gt $d.Paragraphs.Count) {
        $d.Paragraphs($idx - 1).Range.InsertParagraphAfter()
    }
    $p = $d.Paragraphs($idx)
    $p.Range.Text = $texts[$i]
    $p.Format.SpaceAfter = 0
    $p.Format.LineSpacingRule = 0
    if ($indented[$i]) {
        $p.Format.LeftIndent = 36
    }
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
